$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0.8860820675387799, 0.1553318658871028, 0.5503119725020014, 0.1912720824299754, 0, 1.159882681369169, 1.150442336281017, 0, 0.07980958790355608, 0.4232394507347692, 0.3681403120152424, 0, 2.161137300768083, 4.701618427858591),
    @(0.8483674220107389, 0.1551130559053888, 0.5477061812628108, 0.1914511427402736, 0, 1.164778054027352, 1.156711881146251, 0, 0.08021126709107662, 0.3897153628065269, 0.3638353556750502, 0, 2.180654437182186, 4.724838912469522),
    @(0.8255326427194518, 0.1549892857467974, 0.5463355161276269, 0.1916290969168699, 0, 1.168329061258895, 1.160952141431942, 0, 0.08047493473083733, 0.3692265905135486, 0.3613384786793716, 0, 2.193248588815246, 4.741060149123712),
    @(0.8163090747943045, 0.1549415289620946, 0.5458347754806994, 0.191718757068303, 0, 1.169913321686124, 1.1627784692487, 0, 0.08058667354629634, 0.3609017461770634, 0.3603579351156867, 0, 2.198534398690163, 4.748164712210524),
    @(0.8147824712029887, 0.1549337614777997, 0.5457551244670924, 0.1917346816200265, 0, 1.170184676004354, 1.163087675888846, 0, 0.08060548717483407, 0.3595209084257505, 0.3601973527972362, 0, 2.199421381786308, 4.74937428656969),
    @(0.8254079182603959, 0.154988630801018, 0.5463285286584068, 0.1916302366419984, 0, 1.168349871522501, 1.160976373383441, 0, 0.08047642428957058, 0.3691142188100685, 0.3613251049077562, 0, 2.193319253128475, 4.741153961809687),
    @(0.873011646342718, 0.1552542362280676, 0.5493659963404127, 0.1913197308997852, 0, 1.161457498622852, 1.152523058440508, 0, 0.07994455640479536, 0.4116609271894731, 0.3666256658991358, 0, 2.167740129090074, 4.709217551283473),
    @(0.9688896926551251, 0.1558582042605252, 0.557135808113614, 0.1912488132518604, 0, 1.152264303255862, 1.13904072177769, 0, 0.07903634716215002, 0.4958291053474682, 0.3781760311401428, 0, 2.122420427968325, 4.662153054426085),
    @(1.040840820960426, 0.1563515691542037, 0.5639428477043396, 0.1915226612009953, 0, 1.148141772247811, 1.13101436627872, 0, 0.07845072503882378, 0.5580933758871538, 0.3873606555901432, 0, 2.092069194832169, 4.637039601450965),
    @(1.07389488136161, 0.156586575758098, 0.5672766998180805, 0.1917176049738067, 0, 1.146837108268386, 1.127769462727187, 0, 0.07820192598222953, 0.5865072547752277, 0.3916894405079745, 0, 2.07889986995461, 4.627665494497421),
    @(1.086457369287643, 0.1566770677799383, 0.5685731245169734, 0.1918015062297229, 0, 1.146425060601743, 1.126599006081776, 0, 0.07811023505381343, 0.5972792292953955, 0.3933501729978985, 0, 2.074004617455412, 4.624410165023335),
    @(1.083749796489514, 0.156657512305955, 0.5682924080391842, 0.1917829888857554, 0, 1.146510156493818, 1.126848492922917, 0, 0.07812987020410311, 0.5949587566972809, 0.3929915499905832, 0, 2.075054819947627, 4.625098168907186),
    @(1.074927495033165, 0.156593990641781, 0.567382677595333, 0.1917243057782407, 0, 1.146801565824745, 1.127672000411764, 0, 0.07819433197468229, 0.58739322955185, 0.3918256396140549, 0, 2.078495297080013, 4.627391777623728),
    @(1.069529498630089, 0.156555276551515, 0.5668298603503672, 0.1916896722829513, 0, 1.146990739374161, 1.128184013948697, 0, 0.07823414513280014, 0.5827607032578328, 0.3911142834098058, 0, 2.080614628762706, 4.628835013343831),
    @(1.038687091380694, 0.1563364219363805, 0.5637297343289731, 0.1915113338621808, 0, 1.148238513406355, 1.131234595146381, 0, 0.07846733822747964, 0.5562382108944064, 0.3870807781795094, 0, 2.092942671195672, 4.637693444025189),
    @(1.019848442396921, 0.156204855179169, 0.5618885739044828, 0.1914199246544115, 0, 1.149150105010492, 1.133210018331695, 0, 0.07861489795466348, 0.5399900143517584, 0.3846448338705386, 0, 2.100668829427788, 4.643652631716321),
    @(1.00904341423589, 0.1561301777867072, 0.5608519318300011, 0.1913739742671474, 0, 1.149728148903151, 1.134384481416475, 0, 0.07870142776454436, 0.5306529475041941, 0.3832579321870639, 0, 2.105172725443855, 4.647273188406075),
    @(1.005390275044959, 0.156105065043576, 0.5605047850974074, 0.1913595554177441, 0, 1.149933093166666, 1.134788707642329, 0, 0.07873101015215234, 0.5274930501020094, 0.3827907926034726, 0, 2.106707974831457, 4.648532203965203),
    @(1.021850701818352, 0.1562187577156138, 0.5620822571081447, 0.191428969858741, 0, 1.149047505127271, 1.132995773038616, 0, 0.07859901849385231, 0.541718790911176, 0.3849026768605768, 0, 2.099840155473345, 4.642998294961529),
    @(1.077517587011926, 0.1566126079409216, 0.5676489668584566, 0.1917412691616676, 0, 1.146713746896609, 1.127428534452335, 0, 0.07817532955230533, 0.5896150806194669, 0.392167513008161, 0, 2.077482257241399, 4.626710100992284),
    @(1.114164709390195, 0.156878746411202, 0.5714850483257123, 0.1920041124693626, 0, 1.145666430973662, 1.124129892418679, 0, 0.07791313197637528, 0.620989218879231, 0.3970408438955104, 0, 2.063404489994661, 4.617780885856291),
    @(1.094581423926059, 0.1567359107251676, 0.5694196016774811, 0.191858466322703, 0, 1.146181694121879, 1.125859378198896, 0, 0.07805172846214603, 0.604237946697566, 0.394428435794282, 0, 2.070869160404667, 4.622389676531299),
    @(1.020945400647719, 0.1562124693805345, 0.5619946248126837, 0.1914248599547435, 0, 1.149093722454637, 1.133092512573072, 0, 0.0786061923154513, 0.5409371977814033, 0.3847860638117311, 0, 2.100214605762706, 4.643293514598668),
    @(0.9426847999073971, 0.1556860247015024, 0.5548404503920352, 0.1912106123544959, 0, 1.154288894800402, 1.142357524566478, 0, 0.07926766721197076, 0.4729830218427082, 0.3749282069783766, 0, 2.134162987540416, 4.673221496409155)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

